$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.8888888888888888
$ws.Range("N2").Value = 0.8935714285714287

$ws.Range("H3").Value = 0.7777777777777778
$ws.Range("N3").Value = 0.6141666666666666

$ws.Range("H4").Value = 0.6666666666666666
$ws.Range("N4").Value = 0.5753571428571428

$ws.Range("H5").Value = 0.3333333333333333
$ws.Range("N5").Value = 0.4667857142857144

$ws.Range("H6").Value = 0.4444444444444444
$ws.Range("N6").Value = 0.3451190476190477

$ws.Range("H7").Value = 0.3333333333333333
$ws.Range("N7").Value = 0.3097619047619048

$ws.Range("H8").Value = 0.5555555555555556
$ws.Range("N8").Value = 0.2620238095238095
